$d = $word.ActiveDocument

# --- Change 1: collapse the three "3.1" bullet/heading paragraphs (plus the
# trailing blank paragraph) down to a single empty paragraph that keeps only
# the "_1t3h5sf" bookmark. ---

# Locate the paragraph that currently starts with "If no code has been
# written" (the paragraph carrying bookmark _1t3h5sf).
$targetIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "If no code has been written*") {
        $targetIndex = $i
        break
    }
}

# Replace its content with just the bookmark, dropping its pPr
# (Heading2/numPr/spacing/rPr) and all its runs.
$target = $d.Paragraphs.Item($targetIndex)
$bookmarkFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="0" w:name="_1t3h5sf" w:colFirst="0" w:colLast="0"/><w:bookmarkEnd w:id="0"/></w:p>'
$target.Range.InsertXML($bookmarkFrag) | Out-Null

# Remove the following two heading paragraphs ("For a Web application..." /
# "For an application...") together with the blank paragraph that trails
# them, collapsing everything down to the single bookmark-only paragraph.
$pWebApp = $d.Paragraphs.Item($targetIndex + 1)
$pBlank = $d.Paragraphs.Item($targetIndex + 3)
$killRange = $d.Range($pWebApp.Range.Start, $pBlank.Range.End)
$killRange.Delete()

# --- Change 2: merge the "fault" / proofErr-wrapped runs into a single run ---
$old = "if you have a suggestion for how the feature might be implemented in a better way, fault or "
$new = "if you have a suggestion for how the feature might be implemented in a better way, fault or "
$d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
